$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.671.54"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "3.781.95"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'647.94"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "'166.20"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "3.779.96"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'6.93"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").Value = "'34.99"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "4.418.42"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "3.750.00"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "69.566.81"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "'17.77"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "'469.95"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'9.62"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").Value = "'0.711"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'82.10"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("D26").Value = "'12.36"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").Value = "'10.39"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").Value = "'2.12"
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D30").Value = "3.929.10"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'2.28"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").Value = "'7.19"
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("D34").Value = "'28.87"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").Value = "'0.174"
$ws.Range("E35").Value = "  +15.27%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "3.735.64"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").Value = "'8.88"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("D41").Value = "'5.87"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D45").Value = "'45.43"
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "'157.26"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "'47.75"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "'0.300"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "'8.39"
$ws.Range("E51").Value = "  -0.93%  "
